# Fix the typo in the "Keynote: Gil Bronza" entry (-> "Gil Broza") and
# leave the selection on that cell, matching the author's edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 (merged B8:E8) currently holds "Keynote: Gil Bronza" -- correct the name.
$ws.Range("B8").Value = "Keynote: Gil Broza"

# Reflect that B8:E8 is the active selection after the edit.
$ws.Range("B8:E8").Select()
